# Update the "想去人数" (number of people interested) values on both the
# "展览" and "全部类型" worksheets, which hold duplicate data.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Mapping of row -> new value for column F
$updates = @{
    2 = 68
    3 = 1363
    4 = 9
    5 = 6
    7 = 36
    8 = 186
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
